$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure numeric-looking strings (prices, percentages) stay as text,
# matching the original inline-string cell contents (e.g. "1.00", "413.46").
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "62.188.79"
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = "  +0.50%  "

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "3.425.44"
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = "  +0.10%  "

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "0.999"
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = "  -0.06%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "413.46"
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = "  +1.16%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "129.55"
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = "  +0.78%  "

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.624"
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = "  -1.22%  "

$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = "  +0.03%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.723"
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = "  -1.41%  "

$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = "  +1.04%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "42.53"
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = "  -0.18%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "9.28"
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = "  +1.81%  "

$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = "  +6.36%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "3.967.31"
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = "  +0.05%  "

$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = "  -0.36%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "20.52"
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = "  -3.29%  "

$ws.Range("B17").Value = "Uniswap"
$ws.Range("C17").Value = "https://coinranking.com/coin/_H5FVG9iW+uniswap-uni"
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "13.02"
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = "  +5.71%  "

$ws.Range("B18").Value = "WrappedEther"
$ws.Range("C18").Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "3.411.59"
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = "  -1.10%  "

$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = "  -0.61%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "62.205.03"
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = "  +0.55%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "482.56"
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = "  +8.48%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "91.44"
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = "  +0.41%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "3.28"
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = "  +2.72%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "13.42"
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = "  +3.84%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "10.38"
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = "  +18.17%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "3.29"
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = "  +1.72%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "33.36"

$ws.Range("E28").NumberFormat = "@"
$ws.Range("E28").Value = "  -0.39%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "7.60"
$ws.Range("E29").NumberFormat = "@"
$ws.Range("E29").Value = "  +0.04%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "11.94"
$ws.Range("E30").NumberFormat = "@"
$ws.Range("E30").Value = "  -0.11%  "

$ws.Range("E31").NumberFormat = "@"
$ws.Range("E31").Value = "  -4.26%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "0.166"
$ws.Range("E32").NumberFormat = "@"
$ws.Range("E32").Value = "  -3.08%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.111"
$ws.Range("E33").NumberFormat = "@"
$ws.Range("E33").Value = "  -2.86%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "40.48"
$ws.Range("E34").NumberFormat = "@"
$ws.Range("E34").Value = "  -5.23%  "

$ws.Range("E35").NumberFormat = "@"
$ws.Range("E35").Value = "  +0.04%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "58.36"
$ws.Range("E36").NumberFormat = "@"
$ws.Range("E36").Value = "  +9.50%  "

$ws.Range("E37").NumberFormat = "@"
$ws.Range("E37").Value = "  -2.14%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "1.00"
$ws.Range("E38").NumberFormat = "@"
$ws.Range("E38").Value = "  +0.17%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "3.02"
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = "  +3.03%  "

$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = "  +4.04%  "

$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = "  -0.17%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "3.34"
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = "  -1.43%  "

$ws.Range("B43").Value = "Monero"
$ws.Range("C43").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "145.49"
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = "  +2.57%  "

$ws.Range("B44").Value = "WEMIXToken"
$ws.Range("C44").Value = "https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "2.68"
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = "  +7.27%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "2.06"
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = "  +4.29%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "4.29"
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = "  +0.65%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "2.42"
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = "  +15.44%  "

$ws.Range("B48").Value = "PEPE"
$ws.Range("C48").Value = "https://coinranking.com/coin/03WI8NQPF+pepe-pepe"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.0₃0553"
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = "  +35.28%  "

$ws.Range("B49").Value = "Celestia"
$ws.Range("C49").Value = "https://coinranking.com/coin/YQcD0lBl7+celestia-tia"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "16.32"
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = "  -1.41%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "22.36"
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = "  -0.11%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "111.18"
$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = "  +7.81%  "
